# Update the "insert-contracts" sheet so that the CONCATENATE formulas
# build lower-case "insert into ... values (" SQL statements instead of
# the previous upper-case "INSERT INTO ... VALUES (" ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("insert-contracts")

for ($i = 2; $i -le 8; $i++) {
    $formula = '=CONCATENATE("insert into ",''create-contracts''!$B$1," values (","''",raw!A' + $i + ',"'', ","''",raw!B' + $i + ',"'', ",raw!C' + $i + ',", ","''",raw!D' + $i + ',"'', ","''",raw!E' + $i + ',"'', ","''",raw!F' + $i + ',"'', ","''",raw!G' + $i + ',"'');")'
    $ws.Range("A$i").Formula = $formula
}

# Row 2 grows a bit taller to match the new content layout.
$ws.Rows.Item(2).RowHeight = 60
